$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for the "ByStages" columns to include the range qualifiers
$ws.Range("E1").Value = "ByStages nc (10-250)"
$ws.Range("F1").Value = "ByStages high nb (250-4000)"

# Mark the currently running benchmark cases in the "ByStages nc" column
$ws.Range("E5").Value = "Running"
$ws.Range("E6").Value = "Running"
$ws.Range("E11").Value = "Running"

# Give the newly-used C11 cell the same highlight fill used by the rest of column C
$ws.Range("C11").Interior.Color = $ws.Range("C2").Interior.Color

# Widen column E so the new, longer labels fit
$ws.Columns("E").ColumnWidth = 17.83

# Move the active selection to match where the author left off editing
$ws.Range("E10").Select()
